$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.688.11"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "2.199.11"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'258.74"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").Value = "'83.44"
$ws.Range("E6").Value = "  +11.14%  "
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "'44.69"
$ws.Range("E10").Value = "  +11.07%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "'7.21"
$ws.Range("E12").Value = "  +6.44%  "
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").Value = "2.526.59"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "'14.36"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "2.234.81"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "43.615.44"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").Value = "'69.82"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").Value = "'5.92"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("E22").Value = "  +12.25%  "
$ws.Range("D23").Value = "'230.33"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'8.97"
$ws.Range("E24").Value = "  -5.34%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'10.68"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").Value = "'39.40"
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  +3.96%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.21"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("D31").Value = "'174.30"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "'20.41"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  +5.88%  "
$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("E36").Value = "  +4.26%  "
$ws.Range("D37").Value = "'4.54"
$ws.Range("E37").Value = "  +7.95%  "
$ws.Range("D38").Value = "'0.0359"
$ws.Range("E38").Value = "  +7.41%  "
$ws.Range("D39").Value = "'12.57"
$ws.Range("E39").Value = "  +5.21%  "
$ws.Range("D40").Value = "'2.84"
$ws.Range("E40").Value = "  +10.52%  "
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  +7.72%  "
$ws.Range("D43").Value = "'5.49"
$ws.Range("E43").Value = "  +6.26%  "
$ws.Range("D44").Value = "'0.200"
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0979"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.31"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'99.83"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").Value = "'1.19"
$ws.Range("E48").Value = "  +5.65%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").Value = "'0.444"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.11"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "'1.47"
$ws.Range("E51").Value = "  +8.60%  "
